$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.348.70"
$ws.Range("E2").Value = "  -3.42%  "
$ws.Range("D3").Value = "1.648.32"
$ws.Range("E3").Value = "  -3.42%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'213.76"
$ws.Range("E5").Value = "  -1.83%  "
$ws.Range("E6").Value = "  -2.65%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").Value = "'24.02"
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("E9").Value = "  -1.34%  "
$ws.Range("E10").Value = "  -2.47%  "
$ws.Range("D11").Value = "'0.0876"
$ws.Range("E11").Value = "  -1.38%  "
$ws.Range("D12").Value = "1.883.39"
$ws.Range("E12").Value = "  -3.49%  "
$ws.Range("D13").Value = "1.653.54"
$ws.Range("E13").Value = "  -3.17%  "
$ws.Range("E14").Value = "  -2.67%  "
$ws.Range("D15").Value = "'0.566"
$ws.Range("E15").Value = "  +1.15%  "
$ws.Range("D16").Value = "'65.75"
$ws.Range("E16").Value = "  -2.24%  "
$ws.Range("D17").Value = "27.360.51"
$ws.Range("E17").Value = "  -3.24%  "
$ws.Range("D18").Value = "'234.34"
$ws.Range("E18").Value = "  -7.16%  "
$ws.Range("E19").Value = "  -2.83%  "
$ws.Range("D20").Value = "'7.45"
$ws.Range("E20").Value = "  -3.26%  "
$ws.Range("E21").Value = "  +0.18%  "
$ws.Range("D22").Value = "'4.41"
$ws.Range("E22").Value = "  -3.30%  "
$ws.Range("D23").Value = "'9.27"
$ws.Range("E23").Value = "  -3.16%  "
$ws.Range("E24").Value = "  -0.23%  "
$ws.Range("D25").Value = "'145.97"
$ws.Range("E25").Value = "  -1.08%  "
$ws.Range("D26").Value = "'7.16"
$ws.Range("E26").Value = "  -2.59%  "
$ws.Range("D27").Value = "'16.05"
$ws.Range("E27").Value = "  -3.17%  "
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("E29").Value = "  -2.63%  "
$ws.Range("E30").Value = "  -2.69%  "
$ws.Range("E31").Value = "  -1.51%  "
$ws.Range("D33").Value = "1.458.49"
$ws.Range("E34").Value = "  -3.50%  "
$ws.Range("E35").Value = "  -5.16%  "
$ws.Range("E36").Value = "  -0.32%  "
$ws.Range("D37").Value = "'0.908"
$ws.Range("E37").Value = "  -5.53%  "
$ws.Range("E38").Value = "  -4.16%  "
$ws.Range("E39").Value = "  -2.59%  "
$ws.Range("E40").Value = "  -0.34%  "
$ws.Range("E41").Value = "  -0.01%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").Value = "'65.18"
$ws.Range("E42").Value = "  -6.07%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.42"
$ws.Range("E43").Value = "  -3.87%  "
$ws.Range("D44").Value = "'2.21"
$ws.Range("E44").Value = "  -1.83%  "
$ws.Range("D45").Value = "1.791.19"
$ws.Range("E45").Value = "  -3.44%  "
$ws.Range("E46").Value = "  -2.40%  "
$ws.Range("E47").Value = "  -0.44%  "
$ws.Range("D48").Value = "'88.09"
$ws.Range("E48").Value = "  -1.96%  "
$ws.Range("E49").Value = "  -4.78%  "
$ws.Range("E50").Value = "  -2.20%  "
$ws.Range("D51").Value = "'7.75"
$ws.Range("E51").Value = "  -3.44%  "
